# Remove the specific account rows that were dropped from the "Export" sheet.
# (Identified by matching Conta/Nome/Saldo triples against the before/after
# state; these 23 rows are deleted outright, nothing else moves or changes.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# 1-based worksheet row numbers (header is row 1) to delete, computed from the
# original row order. Must be removed highest-to-lowest so earlier deletions
# don't shift the numbering of rows still queued for removal.
$rowsToDelete = @(53,50,49,48,47,46,44,42,41,39,38,36,34,33,32,31,27,26,25,21,19,17,14)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
